$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15 updates: add class/code-review hours and a note about 3 hours code review
$ws.Range("D15").Value = 3
$ws.Range("F15").Value = 8
$ws.Range("I15").Value = "3 hours code review"

# Update the active selection to match the author's last edited cell
$ws.Range("I16").Select()
